$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: add header values for new columns P and Q (continue 0..15 sequence) ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the header style (bold/border/centered) from O1 onto the two new header cells
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Rows 2-25: swap I<->K and M<->O, then fill new columns P and Q with 2 ---
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value()   # column I = 9
    $kVal = $ws.Cells.Item($r, 11).Value()  # column K = 11
    $ws.Cells.Item($r, 9).Value = $kVal
    $ws.Cells.Item($r, 11).Value = $iVal

    $mVal = $ws.Cells.Item($r, 13).Value()  # column M = 13
    $oVal = $ws.Cells.Item($r, 15).Value()  # column O = 15
    $ws.Cells.Item($r, 13).Value = $oVal
    $ws.Cells.Item($r, 15).Value = $mVal

    $ws.Cells.Item($r, 16).Value = 2  # column P = 16
    $ws.Cells.Item($r, 17).Value = 2  # column Q = 17
}
